# Atualizacao dos dados: 27.12.2025 09:38
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13 already carries the exact formatting (borders/number formats) that a
# filled-in data row needs, so copy it straight down into row 14 - which so
# far only held the "id" placeholder in column A - and then overwrite the
# copied values with the new entry (id = 13).
$ws.Range("A13:L13").Copy($ws.Range("A14"))

$ws.Range("A14").Value = 13
$ws.Range("B14").Value = 23
$ws.Range("C14").Value = 69725
$ws.Range("D14").Formula = "=C14*F14"
$ws.Range("E14").Value = 106
$ws.Range("F14").Value = 1.5
$ws.Range("G14").Value = 1588
$ws.Range("H14").Value = 0.050231481481481481
$ws.Range("I14").Value = 7000
$ws.Range("J14").Value = "Vampiro"
$ws.Range("K14").Value = "Desafio"
$ws.Range("L14").Value = 46018

# The remaining placeholder rows (15:27) no longer need their running "id"
# numbers in column A - wipe the cells (not just their contents) so the row
# goes back to holding only the empty, pre-formatted J cell.
$ws.Range("A15:A51").Clear()

# The tail of empty placeholder rows (28:51) is no longer needed at all.
$ws.Range("A28:A51").EntireRow.Delete()

# Leave the selection where the edit session ended.
$ws.Range("M19").Select()
